# Updates the cryptos list: refreshed prices/volumes for most rows, and
# swaps the NEARProtocol / Monero rows (43 <-> 44) to match the new ranking.
# Column D ("Price") cells are forced to Text format before assignment
# (and the format reset afterwards) so that values which look numeric
# (e.g. "127.59") are stored as text, matching the original inlineStr
# string cells rather than being auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.489.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.987.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "382.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0858"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.455.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.998.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.501.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.91%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.027.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.282.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0333"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
